# Apply the "Add files via upload" edit to testData.xlsx:
#   - G2 ("Under 23" answer for the sample test case) changes from "No" to "Yes"
#     (this also introduces a new shared string "Yes" alongside the existing "No").
#   - The sheet's active selection moves from D6 to D3.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("insuranceData")

# Update the data cell: Under 23 -> "Yes"
$ws.Range("G2").Value = "Yes"

# Update the current selection/active cell shown in the saved view.
$ws.Range("D3").Select()
